$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9 (ALC)
$ws.Range("H9").Value = 59.75
$ws.Range("I9").Value = 60.636364
$ws.Range("K9").Value = 60.636364
$ws.Range("M9").Value = 108.363636

# Row 107 (ALC)
$ws.Range("H107").Value = 195.58333
$ws.Range("I107").Value = 122.454544
$ws.Range("K107").Value = 122.454544
$ws.Range("M107").Value = 1797.545456

$ws = $wb.Worksheets.Item("ARM")
# Row 45 (ARM)
$ws.Range("H45").Value = 4783.1665
$ws.Range("I45").Value = 2500
$ws.Range("K45").Value = 2500
$ws.Range("M45").Value = -2123

# Row 74 (ARM)
$ws.Range("H74").Value = 4087.7
$ws.Range("I74").Value = 3044.4707
$ws.Range("J74").Value = 9999.333000000001
$ws.Range("K74").Value = 3044.4707
$ws.Range("L74").Value = 9999.333000000001
$ws.Range("M74").Value = -2170.4707
$ws.Range("N74").Value = -11747.333

# Row 77 (ARM)
$ws.Range("H77").Value = 4087.7
$ws.Range("I77").Value = 3044.4707
$ws.Range("J77").Value = 9999.333000000001
$ws.Range("K77").Value = 15222.3535
$ws.Range("L77").Value = 49996.665
$ws.Range("M77").Value = -10854.3535
$ws.Range("N77").Value = -58732.665

# Row 102 (ARM)
$ws.Range("H102").Value = 4316.4546
$ws.Range("I102").Value = 2185.3125
$ws.Range("K102").Value = 2185.3125
$ws.Range("M102").Value = -563.3125

$ws = $wb.Worksheets.Item("BSM")
# Row 51 (BSM)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

# Row 102 (BSM)
$ws.Range("H102").Value = 11750
$ws.Range("I102").Value = 11750
$ws.Range("K102").Value = 11750
$ws.Range("M102").Value = -8505

# Row 139 (BSM)
$ws.Range("H139").Value = 69999.336
$ws.Range("J139").Value = 69999
$ws.Range("L139").Value = 69999
$ws.Range("N139").Value = -80279

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 7464
$ws.Range("I31").Value = 7464
$ws.Range("K31").Value = 7464
$ws.Range("M31").Value = -7169

# Row 34 (CRP)
$ws.Range("H34").Value = 7464
$ws.Range("I34").Value = 7464
$ws.Range("K34").Value = 7464
$ws.Range("M34").Value = -7262

# Row 62 (CRP)
$ws.Range("H62").Value = 3005
$ws.Range("I62").Value = 3005
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3005
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2381
$ws.Range("N62").ClearContents()

# Row 65 (CRP)
$ws.Range("H65").Value = 3005
$ws.Range("I65").Value = 3005
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 15025
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -11905
$ws.Range("N65").ClearContents()

# Row 99 (CRP)
$ws.Range("H99").Value = 7182.875
$ws.Range("I99").Value = 6858.75
$ws.Range("J99").Value = 7507
$ws.Range("K99").Value = 6858.75
$ws.Range("L99").Value = 7507
$ws.Range("M99").Value = -5360.75
$ws.Range("N99").Value = -10503

# Row 107 (CRP)
$ws.Range("H107").Value = 537.8889
$ws.Range("I107").Value = 514
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 514
$ws.Range("L107").Value = 600
$ws.Range("M107").Value = 1406
$ws.Range("N107").Value = -4440

# Row 112 (CRP)
$ws.Range("H112").Value = 49999
$ws.Range("J112").Value = 49999
$ws.Range("L112").Value = 49999
$ws.Range("N112").Value = -52953

# Row 126 (CRP)
$ws.Range("H126").Value = 7182.875
$ws.Range("I126").Value = 6858.75
$ws.Range("J126").Value = 7507
$ws.Range("K126").Value = 20576.25
$ws.Range("L126").Value = 22521
$ws.Range("M126").Value = -18106.25
$ws.Range("N126").Value = -27461

# Row 134 (CRP)
$ws.Range("H134").Value = 3280.6
$ws.Range("J134").Value = 6969
$ws.Range("L134").Value = 20907
$ws.Range("N134").Value = -25977

$ws = $wb.Worksheets.Item("CUL")
# Row 17 (CUL)
$ws.Range("H17").Value = 351.1111
$ws.Range("J17").Value = 60
$ws.Range("L17").Value = 180
$ws.Range("N17").Value = -518

# Row 55 (CUL)
$ws.Range("H55").Value = 4208
$ws.Range("J55").Value = 6535.5713
$ws.Range("L55").Value = 19606.7139
$ws.Range("N55").Value = -19960.7139

# Row 68 (CUL)
$ws.Range("H68").Value = 572.5714
$ws.Range("J68").Value = 499.5
$ws.Range("L68").Value = 1498.5
$ws.Range("N68").Value = -3120.5

# Row 71 (CUL)
$ws.Range("H71").Value = 572.5714
$ws.Range("J71").Value = 499.5
$ws.Range("L71").Value = 4495.5
$ws.Range("N71").Value = -12607.5

# Row 80 (CUL)
$ws.Range("H80").Value = 4895.8335
$ws.Range("I80").Value = 4777.778
$ws.Range("K80").Value = 14333.334
$ws.Range("M80").Value = -13397.334

# Row 83 (CUL)
$ws.Range("H83").Value = 4895.8335
$ws.Range("I83").Value = 4777.778
$ws.Range("K83").Value = 43000.002
$ws.Range("M83").Value = -38320.002

# Row 122 (CUL)
$ws.Range("H122").Value = 654
$ws.Range("J122").Value = 814
$ws.Range("L122").Value = 7326
$ws.Range("N122").Value = -12226

# Row 126 (CUL)
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 35 (GSM)
$ws.Range("H35").Value = 23762.25
$ws.Range("I35").Value = 9507.5
$ws.Range("J35").Value = 38017
$ws.Range("K35").Value = 9507.5
$ws.Range("L35").Value = 38017
$ws.Range("M35").Value = -9209.5
$ws.Range("N35").Value = -38613

# Row 92 (GSM)
$ws.Range("H92").Value = 7545.2
$ws.Range("J92").Value = 7994.6665
$ws.Range("L92").Value = 7994.6665
$ws.Range("N92").Value = -11738.6665

# Row 122 (GSM)
$ws.Range("H122").Value = 3161.25
$ws.Range("I122").Value = 1988.8334
$ws.Range("J122").Value = 4333.6665
$ws.Range("K122").Value = 5966.5002
$ws.Range("L122").Value = 13000.9995
$ws.Range("M122").Value = -3516.5002
$ws.Range("N122").Value = -17900.9995

# Row 132 (GSM)
$ws.Range("H132").Value = 57164.227
$ws.Range("I132").Value = 71918.53
$ws.Range("J132").Value = 6999.6
$ws.Range("K132").Value = 215755.59
$ws.Range("L132").Value = 20998.8
$ws.Range("M132").Value = -213225.59
$ws.Range("N132").Value = -26058.8

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Range("H61").Value = 4660.087
$ws.Range("I61").Value = 3886.375
$ws.Range("K61").Value = 3886.375
$ws.Range("M61").Value = -3684.375

# Row 100 (LTW)
$ws.Range("H100").Value = 4347.7856
$ws.Range("I100").Value = 3079
$ws.Range("J100").Value = 9000
$ws.Range("K100").Value = 3079
$ws.Range("L100").Value = 9000
$ws.Range("M100").Value = -2538
$ws.Range("N100").Value = -10082

# Row 113 (LTW)
$ws.Range("H113").Value = 4660.087
$ws.Range("I113").Value = 3886.375
$ws.Range("K113").Value = 3886.375
$ws.Range("M113").Value = -1716.375

# Row 122 (LTW)
$ws.Range("H122").Value = 5033.625
$ws.Range("J122").Value = 5115.385
$ws.Range("L122").Value = 15346.155
$ws.Range("N122").Value = -20246.155

# Row 132 (LTW)
$ws.Range("H132").Value = 9336.3125
$ws.Range("I132").Value = 7939.8
$ws.Range("J132").Value = 11663.833
$ws.Range("K132").Value = 23819.4
$ws.Range("L132").Value = 34991.499
$ws.Range("M132").Value = -21289.4
$ws.Range("N132").Value = -40051.499

# Row 136 (LTW)
$ws.Range("H136").Value = 8699.799999999999
$ws.Range("I136").Value = 8374.75
$ws.Range("K136").Value = 25124.25
$ws.Range("M136").Value = -22574.25

$ws = $wb.Worksheets.Item("WVR")
# Row 54 (WVR)
$ws.Range("H54").Value = 51817.41
$ws.Range("J54").Value = 51817.41
$ws.Range("L54").Value = 51817.41
$ws.Range("N54").Value = -52857.41

# Row 62 (WVR)
$ws.Range("H62").Value = 8166.6665
$ws.Range("J62").Value = 9166.666999999999
$ws.Range("L62").Value = 9166.666999999999
$ws.Range("N62").Value = -10414.667

# Row 65 (WVR)
$ws.Range("H65").Value = 8166.6665
$ws.Range("J65").Value = 9166.666999999999
$ws.Range("L65").Value = 45833.335
$ws.Range("N65").Value = -52073.335

# Row 76 (WVR)
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

# Row 79 (WVR)
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

# Row 136 (WVR)
$ws.Range("H136").Value = 9154.4
$ws.Range("J136").Value = 8645.5
$ws.Range("L136").Value = 25936.5
$ws.Range("N136").Value = -31036.5
